$wb = $excel.ActiveWorkbook

# --- New sheet: Authentication (contains login test data) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$auth = $wb.Worksheets.Add($null, $lastSheet)
$auth.Name = "Authentication"

$auth.Cells.Item(2, 1).NumberFormat = "@"
$auth.Cells.Item(2, 1).Value = "adilkhaleque429@gmail.com"
$auth.Cells.Item(2, 2).NumberFormat = "@"
$auth.Cells.Item(2, 2).Value = "Testunbound6A"
$auth.Cells.Item(2, 3).NumberFormat = "@"
$auth.Cells.Item(2, 3).Value = "Adil"

$auth.Hyperlinks.Add($auth.Cells.Item(2, 1), "mailto:adilkhaleque429@gmail.com")

# --- New sheet: ChangeLanguage (contains language test data) ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$lang = $wb.Worksheets.Add($null, $lastSheet2)
$lang.Name = "ChangeLanguage"

$lang.Cells.Item(2, 1).Value = "Español"
